$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in column J (rows 1-15) while keeping formatting/style.
$ws.Range("J1:J15").ClearContents()

# Move the active selection to P13 to reflect the saved cursor position.
$ws.Range("P13").Select()
